# Applies the cryptos list update described by the commit
# "Updated cryptos list on Mon Nov 25 11:11:34 UTC 2024 with GitHub Actions".
# Only the B (Coin), C (Link), D (Price) and E (Volume(1h)) cells that
# actually changed are touched; everything else is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "98.464.85"
$ws.Range("E2").Value = "  +0.92%  "

# Row 3
$ws.Range("D3").Value = "3.496.47"
$ws.Range("E3").Value = "  +4.02%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.93"
$ws.Range("E5").Value = "  +1.44%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "667.48"
$ws.Range("E6").Value = "  +1.19%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.49"
$ws.Range("E7").Value = "  +6.44%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.429"
$ws.Range("E8").Value = "  +2.32%  "

# Row 9
$ws.Range("E9").Value = "  +3.49%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"

# Row 11
$ws.Range("D11").Value = "3.493.82"
$ws.Range("E11").Value = "  +3.93%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.77"
$ws.Range("E12").Value = "  +12.44%  "

# Row 13
$ws.Range("E13").Value = "  +0.65%  "

# Row 14
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.26"
$ws.Range("E14").Value = "  +2.32%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "98.216.40"
$ws.Range("E15").Value = "  +0.96%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000264"
$ws.Range("E16").Value = "  +3.16%  "

# Row 17
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "4.157.03"
$ws.Range("E17").Value = "  +4.49%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.00"
$ws.Range("E18").Value = "  +5.07%  "

# Row 19
$ws.Range("D19").Value = "3.512.68"
$ws.Range("E19").Value = "  +4.30%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.82"
$ws.Range("E20").Value = "  +11.24%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.86"
$ws.Range("E21").Value = "  +9.62%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.527"
$ws.Range("E22").Value = "  -4.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "525.36"
$ws.Range("E23").Value = "  +5.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.46"
$ws.Range("E24").Value = "  +2.98%  "

# Row 25
$ws.Range("E25").Value = "  +1.92%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.81"
$ws.Range("E26").Value = "  +9.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "98.60"
$ws.Range("E27").Value = "  +4.05%  "

# Row 28
$ws.Range("E28").Value = "  +4.96%  "

# Row 29
$ws.Range("D29").Value = "3.688.70"
$ws.Range("E29").Value = "  +3.89%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.49"
$ws.Range("E30").Value = "  +12.14%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.89"
$ws.Range("E31").Value = "  +13.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.146"
$ws.Range("E32").Value = "  -1.11%  "

# Row 33
$ws.Range("E33").Value = "  +0.60%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.191"
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.602"
$ws.Range("E35").Value = "  +8.78%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.23"
$ws.Range("E36").Value = "  +9.39%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.996"
$ws.Range("E37").Value = "  -0.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.53"
$ws.Range("E38").Value = "  +4.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.01"
$ws.Range("E39").Value = "  +3.25%  "

# Row 40
$ws.Range("E40").Value = "  +4.33%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "529.01"
$ws.Range("E41").Value = "  +1.02%  "

# Row 42
$ws.Range("E42").Value = "  +0.12%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.915"
$ws.Range("E43").Value = "  +8.41%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.78"
$ws.Range("E44").Value = "  +7.54%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.42"
$ws.Range("E45").Value = "  -0.97%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0437"
$ws.Range("E46").Value = "  +4.07%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.80"
$ws.Range("E47").Value = "  +3.10%  "

# Row 48
$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.66"
$ws.Range("E48").Value = "  +0.39%  "

# Row 49
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.68"
$ws.Range("E49").Value = "  -1.86%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.24"
$ws.Range("E50").Value = "  +11.36%  "

# Row 51
$ws.Range("E51").Value = "  +2.94%  "
